# Refactor of the SAP -> Excel export routine:
#   * MIGO / MIRO material-document numbers now line up with the correct
#     rows (row indexes were off by one in the old version of the macro).
#   * The MIRO result (column S) was never written before - add it.
#   * Document numbers coming back from SAP are free-form strings (they can
#     carry leading zeros), so they are written out as text, not numbers.
#   * Selection / view state is left where the macro finished (row 4,
#     columns R:S) and those two columns are widened so the new numbers are
#     fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (PO item 6) ----------------------------------------------------
# Purchase order / purchase requisition numbers generated by this MIGO/MIRO
# run (previously pointed at the wrong pair of documents).
$ws.Range("P2").Value = 4600244336
$ws.Range("Q2").Value = 4503342117

# MIGO material-document number for this row - exported as text because the
# SAP document number field is alphanumeric.
$ws.Range("R2").Value = "'5014660856"

# MIRO (invoice) document number - this column used to be skipped entirely.
$ws.Range("S2").Value = "'5600000372"

# --- Row 3 (PO item 7) ----------------------------------------------------
$ws.Range("P3").Value = 4600244337
$ws.Range("Q3").Value = 4503342118

# MIGO document number for the second row.
$ws.Range("R3").Value = "'5014660857"

# Keep the header row's explicit height now that more rows carry data.
$ws.Rows(1).RowHeight = 15.75

# Widen the two new "document number" columns so the values aren't cut off.
$ws.Range("R:S").ColumnWidth = 10.1

# Leave the sheet scrolled/selected where the export loop finished.
$ws.Range("R4").Select()

# Re-stamp the classification footer (center section) emitted by the
# workbook's "Publica" sensitivity label on save.
$ws.PageSetup.CenterFooter = "&`"Calibri`"&10 &K000000`r# Pública"
